# Applies the "Add data for 2022-06-21" update to the carjacking-by-neighborhood
# workbook: renames the report sheet/title from "through June 12" to
# "through June 13" and bumps the affected monthly counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Rename sheet / update report title -----------------------------------
$ws.Name = "Through 2022-06-13"
$ws.Range("B1").Value = "June 2022 (through June 13)"

# --- Updated counts ---------------------------------------------------------
# Each entry: RowNumber, ColumnLetter, NewValue
$updates = @(
    @(2,  "T",  3),
    @(3,  "B",  4),
    @(4,  "B",  4),
    @(5,  "N",  1),
    @(7,  "AL", 1),
    @(9,  "B",  3),
    @(10, "B",  5),
    @(10, "H",  3),
    @(10, "N",  5),
    @(13, "T",  1),
    @(14, "H",  5),
    @(14, "AF", 1),
    @(17, "H",  1),
    @(19, "N",  2),
    @(23, "B",  1),
    @(32, "B",  1),
    @(37, "B",  1),
    @(41, "B",  1),
    @(46, "T",  1),
    @(51, "B",  1),
    @(57, "B",  2),
    @(66, "H",  1),
    @(75, "B",  1),
    @(92, "AR", 1),
    @(95, "B",  1),
    @(96, "Z",  2)
)

foreach ($u in $updates) {
    $row = $u[0]
    $col = $u[1]
    $val = $u[2]
    $ws.Range("$col$row").Value = $val
}
